{"js": "// Remove the superfluous R-script / commentary paragraphs that followed the\n// \"...8809 to 8752...\" sentence in the \"Data acquisition and cleaning\"\n// section, up to (but not including) the \"Full analysis\" heading.\nconst body = context.document.body;\n\n// Anchor on the sentence that stays right before the block being removed.\nconst results = body.search(\"8809 to 8752\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n    throw new Error(\"Anchor text not found\");\n}\n\nconst anchorPara = results.items[0].paragraphs.getFirst();\nawait context.sync();\n\n// Walk forward from the anchor paragraph, deleting every following\n// paragraph that is still part of the superfluous \"Body Text\" block\n// (the leftover R comments / code + the trailing wrap-up sentence).\n// Stop as soon as we reach a paragraph that is no longer styled\n// \"Body Text\" (that's the \"Full analysis\" Heading 2 paragraph).\nlet current = anchorPara.getNext();\ncurrent.load(\"style\");\nawait context.sync();\n\nlet guard = 0;\nwhile (current.style === \"Body Text\" && guard < 100) {\n    guard++;\n    const toDelete = current;\n    current = current.getNext();\n    current.load(\"style\");\n    toDelete.delete();\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Anchor on the sentence that stays right before the superfluous block.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"8809 to 8752\")\nif (-not $found) {\n    throw \"Anchor text not found\"\n}\n$targetStart = $rng.Start\n\n$paras = $d.Paragraphs\n$n = $paras.Count\n\n# Locate the (1-based) paragraph index that contains the anchor text.\n$anchorIndex = -1\nfor ($i = 1; $i -le $n; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Start -le $targetStart -and $p.Range.End -ge $targetStart) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph not found\"\n}\n\n# Walk forward from the anchor paragraph collecting the run of \"Body Text\"\n# styled paragraphs that follow it (the leftover R comments / code plus the\n# trailing wrap-up sentence). Stop as soon as a paragraph with a different\n# style is hit (the \"Full analysis\" Heading 2 paragraph).\n$blockStartIndex = $anchorIndex + 1\n$blockEndIndex = $anchorIndex\n$i = $blockStartIndex\nwhile ($i -le $n) {\n    $p = $paras.Item($i)\n    if ($p.Style.NameLocal -ne \"Body Text\") {\n        break\n    }\n    $blockEndIndex = $i\n    $i = $i + 1\n}\n\nif ($blockEndIndex -ge $blockStartIndex) {\n    $startP = $paras.Item($blockStartIndex)\n    $endP = $paras.Item($blockEndIndex)\n    $delRange = $d.Range($startP.Range.Start, $endP.Range.End)\n    $delRange.Delete()\n}\n"}
